# Add header row to assignments
#
# The assignments sheet lists, for each course outcome, the assignments
# ("Final Exam", "Final Project") used to assess it. This adds a header
# row above the data so the columns are labelled:
#   Course Outcome | Assignment 1 | Assignment 2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push all existing rows down by one to make room for the new header.
$ws.Rows("1:1").Insert()

# Fill in the new header row.
$ws.Range("A1").Value = "Course Outcome"
$ws.Range("B1").Value = "Assignment 1"
$ws.Range("C1").Value = "Assignment 2"

# Keep the previously-selected "first empty outcome row" selection in sync
# with the data having shifted down by one row (was A5, now A6).
$ws.Range("A6").Select()
